$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link (plain text) updates ---
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

# --- Price / Volume columns: force text storage to match inlineStr format ---
$numericTextCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "E27", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50")
foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.43"
$ws.Range("E2").Value = "0.95%"
$ws.Range("D3").Value = "29.37"
$ws.Range("E3").Value = "-2.83%"
$ws.Range("D4").Value = "5.148"
$ws.Range("E4").Value = "0.10%"
$ws.Range("D5").Value = "0.05779"
$ws.Range("E5").Value = "1.92%"
$ws.Range("D6").Value = "6.628"
$ws.Range("E6").Value = "1.62%"
$ws.Range("D7").Value = "0.8563"
$ws.Range("E7").Value = "1.89%"
$ws.Range("D8").Value = "0.8558"
$ws.Range("E8").Value = "-0.89%"
$ws.Range("D9").Value = "0.1365"
$ws.Range("E9").Value = "2.04%"
$ws.Range("D10").Value = "0.07020"
$ws.Range("E10").Value = "1.48%"
$ws.Range("D11").Value = "0.03050"
$ws.Range("E11").Value = "6.49%"
$ws.Range("D12").Value = "0.09371"
$ws.Range("E12").Value = "-0.10%"
$ws.Range("D13").Value = "0.001545"
$ws.Range("E13").Value = "2.28%"
$ws.Range("D14").Value = "0.0006008"
$ws.Range("E14").Value = "0.32%"
$ws.Range("D15").Value = "0.005981"
$ws.Range("E15").Value = "-2.33%"
$ws.Range("D16").Value = "3.486"
$ws.Range("E16").Value = "-0.64%"
$ws.Range("D17").Value = "3.165"
$ws.Range("E17").Value = "4.74%"
$ws.Range("D18").Value = "2.166"
$ws.Range("E18").Value = "1.85%"
$ws.Range("D19").Value = "0.3204"
$ws.Range("E19").Value = "1.68%"
$ws.Range("D20").Value = "0.03320"
$ws.Range("E20").Value = "1.60%"
$ws.Range("E21").Value = "-1.09%"
$ws.Range("D22").Value = "3.175"
$ws.Range("E22").Value = "-12.50%"
$ws.Range("D23").Value = "0.04138"
$ws.Range("E23").Value = "-0.43%"
$ws.Range("D24").Value = "0.1400"
$ws.Range("E24").Value = "1.88%"
$ws.Range("E25").Value = "1.30%"
$ws.Range("E26").Value = "-4.27%"
$ws.Range("E27").Value = "2.58%"
$ws.Range("E28").Value = "3.37%"
$ws.Range("D40").Value = "0.03728"
$ws.Range("E40").Value = "0.38%"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "1.03%"
$ws.Range("D42").Value = "0.002450"
$ws.Range("E42").Value = "6.03%"
$ws.Range("D43").Value = "0.003503"
$ws.Range("E43").Value = "-34.24%"
$ws.Range("D44").Value = "0.008535"
$ws.Range("E44").Value = "-12.56%"
$ws.Range("D45").Value = "0.00005289"
$ws.Range("E45").Value = "3.75%"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("D47").Value = "0.05798"
$ws.Range("E47").Value = "-42.01%"
$ws.Range("D48").Value = "0.002172"
$ws.Range("E48").Value = "-20.08%"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.02%"

foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
